# "Update new sheet log from #50"
# Append one new log row (row 20) to the activity-log worksheet, in the
# same shape as the existing rows (Source Branch, Author, Action,
# Comment, Date, Change ID).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20
$ws.Cells.Item($row, 1).Value = "edit1"
$ws.Cells.Item($row, 2).Value = "riya-morankar"
$ws.Cells.Item($row, 3).Value = "Merged"
$ws.Cells.Item($row, 4).Value = "N/A"

# Column E holds a plain text date like "2025-06-18" in every other row
# (not a real date value). Force text so Excel's date auto-detection
# doesn't convert it to a serial date, then drop the now-unneeded
# number format so the cell ends up with the default (unstyled) look,
# matching the rest of the sheet.
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "2025-06-19"
$ws.Cells.Item($row, 5).ClearFormats()

$ws.Cells.Item($row, 6).Value = "6bd9de3f8a1122122b8fdaf8a1dea58b3a601eed"
